# Auto-generated Excel COM-interop script to apply BRVM automatic update
# (GitHub Actions daily refresh of data/recommandations.xlsx)

$wb = $excel.ActiveWorkbook

# --- Sheet "Recommandations" ---
$ws1 = $wb.Worksheets.Item("Recommandations")

$ws1.Range("D2").Value = 2494.12
$ws1.Range("E2").Value = 111.35

$ws1.Range("D3").Value = 2025
$ws1.Range("E3").Value = 680

$ws1.Range("D4").Value = 2015
$ws1.Range("E4").Value = 680

$ws1.Range("D5").Value = 1988.01
$ws1.Range("E5").Value = 662.89

$ws1.Range("A6").Value = "NEI-CEDA CI"
$ws1.Range("D6").Value = 1780

$ws1.Range("A7").Value = "SETAO CI"
$ws1.Range("D7").Value = 1770
$ws1.Range("E7").Value = 585

$ws1.Range("D8").Value = 1710
$ws1.Range("E8").Value = 585

$ws1.Range("D9").Value = 1580

$ws1.Range("D10").Value = 1074.59
$ws1.Range("E10").Value = 354.23

$ws1.Range("D11").Value = 1042.73
$ws1.Range("E11").Value = 345.13

$ws1.Range("D12").Value = 924.71
$ws1.Range("E12").Value = 308.14

$ws1.Range("D13").Value = 654.65
$ws1.Range("E13").Value = 223.91

$ws1.Range("D14").Value = 553.32
$ws1.Range("E14").Value = 188.09

$ws1.Range("D15").Value = 538.43
$ws1.Range("E15").Value = 180.37

$ws1.Range("D16").Value = 399.42
$ws1.Range("E16").Value = 134.61

$ws1.Range("D17").Value = 388.56
$ws1.Range("E17").Value = 130.29

$ws1.Range("D18").Value = 366.44
$ws1.Range("E18").Value = 122.66

$ws1.Range("D19").Value = 360.14
$ws1.Range("E19").Value = 120.55

$ws1.Range("A20").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Range("D20").Value = 319
$ws1.Range("E20").Value = 106.78

$ws1.Range("A21").Value = "BRVM - ENERGIE"
$ws1.Range("D21").Value = 316.19
$ws1.Range("E21").Value = 103.6

$ws1.Range("D22").Value = 281.37
$ws1.Range("E22").Value = 93.83

$ws1.Range("D23").Value = 22.47
$ws1.Range("E23").Value = 7.49

$ws1.Range("A24").Value = "CIE CI (CIEC)"
$ws1.Range("B24").Value = 1
$ws1.Range("D24").Value = 7.33
$ws1.Range("E24").Value = 7.33

$ws1.Range("A25").Value = "SODE CI (SDCC)"
$ws1.Range("B25").Value = 1
$ws1.Range("D25").Value = 7.26
$ws1.Range("E25").Value = 7.26

$ws1.Range("A26").Value = "SITAB CI (STBC)"
$ws1.Range("D26").Value = 6.35
$ws1.Range("E26").Value = 6.35

$ws1.Range("A27").Value = "UNIWAX CI (UNXC)"
$ws1.Range("D27").Value = 4.46
$ws1.Range("E27").Value = 4.46

$ws1.Range("A28").Value = "FILTISAC CI (FTSC)"
$ws1.Range("D28").Value = 4.33
$ws1.Range("E28").Value = 4.33

$ws1.Range("A29").Value = "SMB CI (SMBC)"
$ws1.Range("D29").Value = 3.92
$ws1.Range("E29").Value = 3.92

$ws1.Range("A30").Value = "BANK OF AFRICA NG (BOAN)"
$ws1.Range("D30").Value = 1.94
$ws1.Range("E30").Value = 5.1

$ws1.Range("A32").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Range("D32").Value = 0.37
$ws1.Range("E32").Value = -5.88

$ws1.Range("A34").Value = "ONATEL BF (ONTBF)"
$ws1.Range("B34").Value = 1
$ws1.Range("D34").Value = -0.28
$ws1.Range("E34").Value = -3.8
$ws1.Range("G34").Value = "👀 À surveiller"

$ws1.Range("A35").Value = "NESTLE CI (NTLC)"
$ws1.Range("D35").Value = -0.85
$ws1.Range("E35").Value = -0.85

$ws1.Range("A36").Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Range("B36").Value = 1
$ws1.Range("D36").Value = -1.47
$ws1.Range("E36").Value = 4.46
$ws1.Range("G36").Value = "👀 À surveiller"

$ws1.Range("A37").Value = "BANK OF AFRICA ML (BOAM)"
$ws1.Range("D37").Value = -1.91
$ws1.Range("E37").Value = -1.91

$ws1.Range("A38").Value = "ORANGE COTE D'IVOIRE (ORAC)"
$ws1.Range("D38").Value = -2.07
$ws1.Range("E38").Value = -2.07

$ws1.Range("A39").Value = "VIVO ENERGY CI (SHEC)"
$ws1.Range("D39").Value = -2.4
$ws1.Range("E39").Value = -2.4

$ws1.Range("A40").Value = "AIR LIQUIDE CI (SIVC)"
$ws1.Range("D40").Value = -2.86
$ws1.Range("E40").Value = -2.86

$ws1.Range("A41").Value = "BERNABE CI (BNBC)"
$ws1.Range("B41").Value = 1
$ws1.Range("D41").Value = -3.95
$ws1.Range("E41").Value = -7.41
$ws1.Range("G41").Value = "👀 À surveiller"

$ws1.Range("A42").Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Range("D42").Value = -4.03
$ws1.Range("E42").Value = -4.03

$ws1.Range("A43").Value = "SETAO CI (STAC)"
$ws1.Range("B43").Value = 0
$ws1.Range("D43").Value = -4.2
$ws1.Range("E43").Value = -4.2
$ws1.Range("G43").Value = "➖ Neutre"

$ws1.Range("A44").Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws1.Range("C44").Value = 2
$ws1.Range("D44").Value = -7.04
$ws1.Range("E44").Value = -4.47

# Row 45 no longer present after the refresh -> remove it (sheet now ends at row 44)
$ws1.Rows.Item(45).Delete()

# --- Sheet "Top_YTD" ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

$ws2.Range("B2").Value = 516539.04

$ws2.Range("B3").Value = 46442.6

$ws2.Range("B4").Value = 45822.5

$ws2.Range("B5").Value = 44261.84

$ws2.Range("A6").Value = "NEI-CEDA CI"
$ws2.Range("B6").Value = 33227

$ws2.Range("A7").Value = "SETAO CI"
$ws2.Range("B7").Value = 32749.18

$ws2.Range("B8").Value = 29964.65

$ws2.Range("B9").Value = 24507.8

$ws2.Range("B10").Value = 9518.48

$ws2.Range("B11").Value = 8865.87

